# Add the "31/12/2022" data point for each of the four variable groups on
# Sheet1. Each group currently occupies 6 rows (years 2012,2017-2021); a new
# row for 2022 is inserted right after each group, pushing the following
# groups down. Insert from the bottom group upward so earlier insert points
# (lower row numbers) are unaffected by later ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) After "Não estuda e não trabalha" block (rows 20-25) -> insert at row 26
$ws.Rows("26:26").Insert()
$ws.Range("A26").Value = "Sergipe"
$ws.Range("B26").Value = "Não estuda e não trabalha"
$ws.Range("C26").Value = "31/12/2022"
$ws.Range("D26").Value = 0.2718082179568723

# 2) After "Só trabalha" block (rows 14-19) -> insert at row 20
$ws.Rows("20:20").Insert()
$ws.Range("A20").Value = "Sergipe"
$ws.Range("B20").Value = "Só trabalha"
$ws.Range("C20").Value = "31/12/2022"
$ws.Range("D20").Value = 0.3535236730639497

# 3) After "Estuda e trabalha" block (rows 8-13) -> insert at row 14
$ws.Rows("14:14").Insert()
$ws.Range("A14").Value = "Sergipe"
$ws.Range("B14").Value = "Estuda e trabalha"
$ws.Range("C14").Value = "31/12/2022"
$ws.Range("D14").Value = 0.08845564513098542

# 4) After "Só estuda" block (rows 2-7) -> insert at row 8
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "Sergipe"
$ws.Range("B8").Value = "Só estuda"
$ws.Range("C8").Value = "31/12/2022"
$ws.Range("D8").Value = 0.2862124638481935

Write-Output "applied edits"
